# "Account TC 01 02 03"
# The test-case sheet TC_01_02 was reopened/edited: cell A2 (previously
# blank but formatted with the quote-prefix style used elsewhere on this
# sheet) was given the literal value "*", and the cursor was left on C2.
# Activating this sheet also makes it the workbook's active tab (and
# correspondingly drops the tabSelected flag that used to sit on
# PriceSheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_01_02")
$ws.Activate()

# Leading apostrophe keeps Excel's "text that looks like it needs a
# quote prefix" handling, so the existing quotePrefix cell style (s="2")
# on A2 is preserved instead of being reset to the default style.
$ws.Range("A2").Value = "'*"

$ws.Range("C2").Select()
